# Weekly update: insert a new price row for "Ajo" (Femacal de La Calera)
# at row 711, pushing the existing rows 711-742 down to 712-743.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 711 (shifts 711..742 -> 712..743)
$ws.Rows.Item(711).Insert()

# Populate the newly inserted row with the new week's data point
$ws.Range("A711").Value = 3
$ws.Range("B711").Value = "Femacal de La Calera"
$ws.Range("C711").Value = "Coquimbo"
$ws.Range("D711").Value = 45041
$ws.Range("E711").Value = 5
$ws.Range("F711").Value = 100112003
$ws.Range("G711").Value = "Ajo"
$ws.Range("H711").Value = "Chino"
$ws.Range("I711").Value = "Primera"
$ws.Range("J711").Value = 85
$ws.Range("K711").Value = 15500
$ws.Range("L711").Value = 16000
$ws.Range("M711").Value = 15765
$ws.Range("N711").Value = "`$/caja 10 kilos"
$ws.Range("O711").Value = "China"
$ws.Range("P711").Value = 1576
$ws.Range("Q711").Value = 10
$ws.Range("R711").Value = "Hortaliza"
